# "same old scenario, make dynamic writing of load shedders."
# Update the "Coupling Parameters" sheet: flip two boolean switches to FALSE
# and change the profit-based dismantling tick value, then move the active
# selection to match where the user ended up (B28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# realistic_candidate_capacities_tobe_installed -> FALSE
$ws.Range("B19").Value = $false

# realistic_candidate_capacities_to_test -> FALSE
$ws.Range("B20").Value = $false

# start_profit_based_dismantling_tick -> 3
$ws.Range("B23").Value = 3

# Reflect the final selection left on the sheet
$ws.Activate()
$ws.Range("B28").Select()
